$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 197.5433703333333
$ws.Range("H2").Value = 592.6301109999999
$ws.Range("I2").Value = 0.3388703761585983
$ws.Range("J2").Value = 0.3388703761585982
$ws.Range("M2").Value = 1.701929666666667
$ws.Range("N2").Value = 5.105789
$ws.Range("O2").Value = 0.02105622887134972
$ws.Range("P2").Value = 0.02105622887134972
$ws.Range("Q2").Value = 336.2049224236198
$ws.Range("R2").Value = 3025.844301812579
$ws.Range("S2").Value = 0.007135332198115816
$ws.Range("T2").Value = 0.007135332198115815

$ws.Range("G3").Value = 197.5433703333333
$ws.Range("H3").Value = 592.6301109999999
$ws.Range("I3").Value = 0.3388703761585983
$ws.Range("J3").Value = 0.3388703761585982
$ws.Range("O3").Value = 0.7732971809418951
$ws.Range("P3").Value = 0.7732971809418953
$ws.Range("Q3").Value = 12347.24034951604
$ws.Range("R3").Value = 111125.1631456443
$ws.Range("S3").Value = 0.2620475065881636
$ws.Range("T3").Value = 0.2620475065881636

$ws.Range("G4").Value = 197.5433703333333
$ws.Range("H4").Value = 592.6301109999999
$ws.Range("I4").Value = 0.3388703761585983
$ws.Range("J4").Value = 0.3388703761585982
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4338690000000001
$ws.Range("N4").Value = 1.301607
$ws.Range("O4").Value = 0.005367815805265532
$ws.Range("P4").Value = 0.005367815805265533
$ws.Range("Q4").Value = 85.707944543153
$ws.Range("R4").Value = 771.371500888377
$ws.Range("S4").Value = 0.0018189937610804
$ws.Range("T4").Value = 0.0018189937610804

$ws.Range("G5").Value = 197.5433703333333
$ws.Range("H5").Value = 592.6301109999999
$ws.Range("I5").Value = 0.3388703761585983
$ws.Range("J5").Value = 0.3388703761585982
$ws.Range("M5").Value = 15.972384
$ws.Range("N5").Value = 47.917152
$ws.Range("O5").Value = 0.1976099128607259
$ws.Range("P5").Value = 0.1976099128607259
$ws.Range("Q5").Value = 3155.238567618208
$ws.Range("R5").Value = 28397.14710856387
$ws.Range("S5").Value = 0.06696414550378202
$ws.Range("T5").Value = 0.06696414550378202

$ws.Range("G6").Value = 197.5433703333333
$ws.Range("H6").Value = 592.6301109999999
$ws.Range("I6").Value = 0.3388703761585983
$ws.Range("J6").Value = 0.3388703761585982
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2157183333333333
$ws.Range("N6").Value = 0.647155
$ws.Range("O6").Value = 0.002668861520763652
$ws.Range("P6").Value = 0.002668861520763652
$ws.Range("Q6").Value = 42.61372660935611
$ws.Range("R6").Value = 383.523539484205
$ws.Range("S6").Value = 0.0009043981074563876
$ws.Range("T6").Value = 0.0009043981074563874

$ws.Range("I7").Value = 0.1369374790620155
$ws.Range("J7").Value = 0.1369374790620154
$ws.Range("M7").Value = 1.701929666666667
$ws.Range("N7").Value = 5.105789
$ws.Range("O7").Value = 0.02105622887134972
$ws.Range("P7").Value = 0.02105622887134972
$ws.Range("Q7").Value = 135.8603695218958
$ws.Range("R7").Value = 1222.743325697062
$ws.Range("S7").Value = 0.002883386900195458
$ws.Range("T7").Value = 0.002883386900195457

$ws.Range("I8").Value = 0.1369374790620155
$ws.Range("J8").Value = 0.1369374790620154
$ws.Range("O8").Value = 0.7732971809418951
$ws.Range("P8").Value = 0.7732971809418953
$ws.Range("S8").Value = 0.1058933665239464
$ws.Range("T8").Value = 0.1058933665239463

$ws.Range("I9").Value = 0.1369374790620155
$ws.Range("J9").Value = 0.1369374790620154
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4338690000000001
$ws.Range("N9").Value = 1.301607
$ws.Range("O9").Value = 0.005367815805265532
$ws.Range("P9").Value = 0.005367815805265533
$ws.Range("Q9").Value = 34.63457028723401
$ws.Range("R9").Value = 311.711132585106
$ws.Range("S9").Value = 0.0007350551644423045
$ws.Range("T9").Value = 0.0007350551644423045

$ws.Range("I10").Value = 0.1369374790620155
$ws.Range("J10").Value = 0.1369374790620154
$ws.Range("M10").Value = 15.972384
$ws.Range("N10").Value = 47.917152
$ws.Range("O10").Value = 0.1976099128607259
$ws.Range("P10").Value = 0.1976099128607259
$ws.Range("Q10").Value = 1275.031533257024
$ws.Range("R10").Value = 11475.28379931322
$ws.Range("S10").Value = 0.02706020330481236
$ws.Range("T10").Value = 0.02706020330481235

$ws.Range("I11").Value = 0.1369374790620155
$ws.Range("J11").Value = 0.1369374790620154
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.2157183333333333
$ws.Range("N11").Value = 0.647155
$ws.Range("O11").Value = 0.002668861520763652
$ws.Range("P11").Value = 0.002668861520763652
$ws.Range("Q11").Value = 17.22020189983222
$ws.Range("R11").Value = 154.98181709849
$ws.Range("S11").Value = 0.0003654671686189915
$ws.Range("T11").Value = 0.0003654671686189914

$ws.Range("G12").Value = 148.824417
$ws.Range("H12").Value = 446.473251
$ws.Range("I12").Value = 0.2552967790580629
$ws.Range("J12").Value = 0.2552967790580629
$ws.Range("M12").Value = 1.701929666666667
$ws.Range("N12").Value = 5.105789
$ws.Range("O12").Value = 0.02105622887134972
$ws.Range("P12").Value = 0.02105622887134972
$ws.Range("Q12").Value = 253.288690416671
$ws.Range("R12").Value = 2279.598213750039
$ws.Range("S12").Value = 0.005375587409964974
$ws.Range("T12").Value = 0.005375587409964974

$ws.Range("G13").Value = 148.824417
$ws.Range("H13").Value = 446.473251
$ws.Range("I13").Value = 0.2552967790580629
$ws.Range("J13").Value = 0.2552967790580629
$ws.Range("O13").Value = 0.7732971809418951
$ws.Range("P13").Value = 0.7732971809418953
$ws.Range("Q13").Value = 9302.1134724739
$ws.Range("R13").Value = 83719.0212522651
$ws.Range("S13").Value = 0.1974202795491459
$ws.Range("T13").Value = 0.1974202795491459

$ws.Range("G14").Value = 148.824417
$ws.Range("H14").Value = 446.473251
$ws.Range("I14").Value = 0.2552967790580629
$ws.Range("J14").Value = 0.2552967790580629
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.4338690000000001
$ws.Range("N14").Value = 1.301607
$ws.Range("O14").Value = 0.005367815805265532
$ws.Range("P14").Value = 0.005367815805265533
$ws.Range("Q14").Value = 64.57030097937302
$ws.Range("R14").Value = 581.1327088143571
$ws.Range("S14").Value = 0.001370386085661252
$ws.Range("T14").Value = 0.001370386085661253

$ws.Range("G15").Value = 148.824417
$ws.Range("H15").Value = 446.473251
$ws.Range("I15").Value = 0.2552967790580629
$ws.Range("J15").Value = 0.2552967790580629
$ws.Range("M15").Value = 15.972384
$ws.Range("N15").Value = 47.917152
$ws.Range("O15").Value = 0.1976099128607259
$ws.Range("P15").Value = 0.1976099128607259
$ws.Range("Q15").Value = 2377.080736900128
$ws.Range("R15").Value = 21393.72663210115
$ws.Range("S15").Value = 0.05044917426328779
$ws.Range("T15").Value = 0.0504491742632878

$ws.Range("G16").Value = 148.824417
$ws.Range("H16").Value = 446.473251
$ws.Range("I16").Value = 0.2552967790580629
$ws.Range("J16").Value = 0.2552967790580629
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.2157183333333333
$ws.Range("N16").Value = 0.647155
$ws.Range("O16").Value = 0.002668861520763652
$ws.Range("P16").Value = 0.002668861520763652
$ws.Range("Q16").Value = 32.104155194545
$ws.Range("R16").Value = 288.937396750905
$ws.Range("S16").Value = 0.0006813517500029639
$ws.Range("T16").Value = 0.0006813517500029639

$ws.Range("G17").Value = 35.426853
$ws.Range("H17").Value = 106.280559
$ws.Range("I17").Value = 0.06077202683121193
$ws.Range("J17").Value = 0.06077202683121192
$ws.Range("M17").Value = 1.701929666666667
$ws.Range("N17").Value = 5.105789
$ws.Range("O17").Value = 0.02105622887134972
$ws.Range("P17").Value = 0.02105622887134972
$ws.Range("Q17").Value = 60.294012117339
$ws.Range("R17").Value = 542.6461090560509
$ws.Range("S17").Value = 0.001279629705933804
$ws.Range("T17").Value = 0.001279629705933804

$ws.Range("G18").Value = 35.426853
$ws.Range("H18").Value = 106.280559
$ws.Range("I18").Value = 0.06077202683121193
$ws.Range("J18").Value = 0.06077202683121192
$ws.Range("O18").Value = 0.7732971809418951
$ws.Range("P18").Value = 0.7732971809418953
$ws.Range("Q18").Value = 2214.318142288791
$ws.Range("R18").Value = 19928.86328059912
$ws.Range("S18").Value = 0.04699483702870139
$ws.Range("T18").Value = 0.04699483702870139

$ws.Range("G19").Value = 35.426853
$ws.Range("H19").Value = 106.280559
$ws.Range("I19").Value = 0.06077202683121193
$ws.Range("J19").Value = 0.06077202683121192
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.4338690000000001
$ws.Range("N19").Value = 1.301607
$ws.Range("O19").Value = 0.005367815805265532
$ws.Range("P19").Value = 0.005367815805265533
$ws.Range("Q19").Value = 15.370613284257
$ws.Range("R19").Value = 138.335519558313
$ws.Range("S19").Value = 0.0003262130461426004
$ws.Range("T19").Value = 0.0003262130461426004

$ws.Range("G20").Value = 35.426853
$ws.Range("H20").Value = 106.280559
$ws.Range("I20").Value = 0.06077202683121193
$ws.Range("J20").Value = 0.06077202683121192
$ws.Range("M20").Value = 15.972384
$ws.Range("N20").Value = 47.917152
$ws.Range("O20").Value = 0.1976099128607259
$ws.Range("P20").Value = 0.1976099128607259
$ws.Range("Q20").Value = 565.851300027552
$ws.Range("R20").Value = 5092.661700247968
$ws.Range("S20").Value = 0.01200915492648548
$ws.Range("T20").Value = 0.01200915492648548

$ws.Range("G21").Value = 35.426853
$ws.Range("H21").Value = 106.280559
$ws.Range("I21").Value = 0.06077202683121193
$ws.Range("J21").Value = 0.06077202683121192
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.2157183333333333
$ws.Range("N21").Value = 0.647155
$ws.Range("O21").Value = 0.002668861520763652
$ws.Range("P21").Value = 0.002668861520763652
$ws.Range("Q21").Value = 7.642221684405
$ws.Range("R21").Value = 68.779995159645
$ws.Range("S21").Value = 0.0001621921239486378
$ws.Range("T21").Value = 0.0001621921239486377

$ws.Range("G22").Value = 121.3248153333333
$ws.Range("H22").Value = 363.974446
$ws.Range("I22").Value = 0.2081233388901116
$ws.Range("J22").Value = 0.2081233388901115
$ws.Range("M22").Value = 1.701929666666667
$ws.Range("N22").Value = 5.105789
$ws.Range("O22").Value = 0.02105622887134972
$ws.Range("P22").Value = 0.02105622887134972
$ws.Range("Q22").Value = 206.4863025186549
$ws.Range("R22").Value = 1858.376722667894
$ws.Range("S22").Value = 0.004382292657139669
$ws.Range("T22").Value = 0.004382292657139668

$ws.Range("G23").Value = 121.3248153333333
$ws.Range("H23").Value = 363.974446
$ws.Range("I23").Value = 0.2081233388901116
$ws.Range("J23").Value = 0.2081233388901115
$ws.Range("O23").Value = 0.7732971809418951
$ws.Range("P23").Value = 0.7732971809418953
$ws.Range("Q23").Value = 7583.279827379454
$ws.Range("R23").Value = 68249.51844641508
$ws.Range("S23").Value = 0.160941191251938
$ws.Range("T23").Value = 0.160941191251938

$ws.Range("G24").Value = 121.3248153333333
$ws.Range("H24").Value = 363.974446
$ws.Range("I24").Value = 0.2081233388901116
$ws.Range("J24").Value = 0.2081233388901115
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 0.4338690000000001
$ws.Range("N24").Value = 1.301607
$ws.Range("O24").Value = 0.005367815805265532
$ws.Range("P24").Value = 0.005367815805265533
$ws.Range("Q24").Value = 52.639076303858
$ws.Range("R24").Value = 473.7516867347221
$ws.Range("S24").Value = 0.001117167747938975
$ws.Range("T24").Value = 0.001117167747938975

$ws.Range("G25").Value = 121.3248153333333
$ws.Range("H25").Value = 363.974446
$ws.Range("I25").Value = 0.2081233388901116
$ws.Range("J25").Value = 0.2081233388901115
$ws.Range("M25").Value = 15.972384
$ws.Range("N25").Value = 47.917152
$ws.Range("O25").Value = 0.1976099128607259
$ws.Range("P25").Value = 0.1976099128607259
$ws.Range("Q25").Value = 1937.846539233088
$ws.Range("R25").Value = 17440.61885309779
$ws.Range("S25").Value = 0.04112723486235827
$ws.Range("T25").Value = 0.04112723486235827

$ws.Range("G26").Value = 121.3248153333333
$ws.Range("H26").Value = 363.974446
$ws.Range("I26").Value = 0.2081233388901116
$ws.Range("J26").Value = 0.2081233388901115
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 0.2157183333333333
$ws.Range("N26").Value = 0.647155
$ws.Range("O26").Value = 0.002668861520763652
$ws.Range("P26").Value = 0.002668861520763652
$ws.Range("Q26").Value = 26.17198695568111
$ws.Range("R26").Value = 235.54788260113
$ws.Range("S26").Value = 0.0005554523707366722
$ws.Range("T26").Value = 0.0005554523707366721

